$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.148.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.83%  '
$ws.Range("D3").Value = "'1.784.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.57%  '
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").Value = "'336.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.66%  '
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("D7").Value = "'0.3824"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.43%  '
$ws.Range("E8").Value = '  -2.41%  '
$ws.Range("D9").Value = "'47.99"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.63%  '
$ws.Range("D10").Value = "'1.188"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.53%  '
$ws.Range("D11").Value = "'0.07452"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.57%  '
$ws.Range("D12").Value = "'1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("D13").Value = "'21.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.29%  '
$ws.Range("D14").Value = "'6.440"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.86%  '
$ws.Range("D15").Value = "'1.780.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.57%  '
$ws.Range("D16").Value = "'7.095"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.19%  '
$ws.Range("D17").Value = "'0.00001091"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.42%  '
$ws.Range("D18").Value = "'0.06653"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.97%  '
$ws.Range("D19").Value = "'83.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.86%  '
$ws.Range("D20").Value = "'1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("D21").Value = "'6.516"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.45%  '
$ws.Range("D22").Value = "'17.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.12%  '
$ws.Range("D23").Value = "'27.135.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.94%  '
$ws.Range("D24").Value = "'12.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -8.16%  '
$ws.Range("D25").Value = "'2.375"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.78%  '
$ws.Range("D26").Value = "'21.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.09%  '
$ws.Range("D27").Value = "'2.495"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.56%  '
$ws.Range("D28").Value = "'1.442"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.19%  '
$ws.Range("D29").Value = "'155.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.09%  '
$ws.Range("D30").Value = "'1.984.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.20%  '
$ws.Range("D31").Value = "'134.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.39%  '
$ws.Range("D32").Value = "'3.989"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.09%  '
$ws.Range("D33").Value = "'6.038"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.40%  '
$ws.Range("D34").Value = "'0.08663"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.80%  '
$ws.Range("D35").Value = "'13.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.13%  '
$ws.Range("D36").Value = "'1.624"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.83%  '
$ws.Range("D37").Value = "'0.6824"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.22%  '
$ws.Range("D38").Value = "'5.394"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.00%  '
$ws.Range("D39").Value = "'0.06295"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.97%  '
$ws.Range("D40").Value = "'0.2175"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.46%  '
$ws.Range("D41").Value = "'0.02322"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.95%  '
$ws.Range("D42").Value = "'1.242"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.12%  '
$ws.Range("D43").Value = "'8.371"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.60%  '
$ws.Range("D44").Value = "'14.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.14%  '
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").Value = "'0.9998"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.17%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = "'0.6418"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.60%  '
$ws.Range("D47").Value = "'3.854"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.13%  '
$ws.Range("D48").Value = "'2.128"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.43%  '
$ws.Range("D49").Value = "'131.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.05%  '
$ws.Range("D50").Value = "'0.07102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.16%  '
$ws.Range("E51").Value = '  -2.25%  '
